# "Generate Report for Handback"
#
# This regenerates the handback timestamps for the file
# 81417268-5735-44f5-abbf-1878c8fbd91c.md across the Overview, zh-cn and
# de-de sheets: Correspond Handoff Datetime / Correspond Handback DateTime
# (and the Overview "Latest HO Xliff Generate Date" roll-up) are refreshed
# to the new handback run's timestamps. The 91930a09-... row is already
# up to date and is left untouched.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the 81417268 file (row 2)
$overview.Range("G2").Value = "2016-08-19 20:52:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 81417268 file (row 2)
$zhcn.Range("H2").Value = "2016-08-19 20:52:15"
$zhcn.Range("K2").Value = "2016-08-19 20:52:32"

# de-de sheet: Correspond Handback DateTime for the 81417268 file (row 2)
$dede.Range("K2").Value = "2016-08-19 20:52:39"
